$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Refresh "last updated" dates for several screens to 12/9/2024 (serial 45635)
# Row 3: Consumer / Gaming
$ws.Range("D3").Value = 45635
# Row 8: Healthcare / Services
$ws.Range("D8").Value = 45635
# Row 11: Tech / Hardware & Semis
$ws.Range("D11").Value = 45635
# Row 12: Tech / Software
$ws.Range("D12").Value = 45635
# Row 19: Hedge Funds
$ws.Range("C19").Value = 45635

# Leave selection on D19, matching the last cell touched in the saved file
$ws.Range("D19").Select()
